$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "26.299.81"
$ws.Range("E2").Value = "  +1.30%  "

$ws.Range("D3").Value = "1.619.95"
$ws.Range("E3").Value = "  +2.01%  "

$ws.Range("E4").Value = "  -0.02%  "

$c = $ws.Range("D5")
$c.Value = "'212.09"
$c.Style = "Normal"

$ws.Range("E6").Value = "  -0.05%  "

$c = $ws.Range("D7")
$c.Value = "'0.484"
$c.Style = "Normal"
$ws.Range("E7").Value = "  +0.96%  "

$ws.Range("E9").Value = "  +0.78%  "

$c = $ws.Range("D10")
$c.Value = "'18.77"
$c.Style = "Normal"

$c = $ws.Range("D11")
$c.Value = "'0.0816"
$c.Style = "Normal"
$ws.Range("E11").Value = "  +0.80%  "

$ws.Range("D12").Value = "1.845.26"
$ws.Range("E12").Value = "  +1.98%  "

$ws.Range("D13").Value = "1.618.76"
$ws.Range("E13").Value = "  +1.89%  "

$c = $ws.Range("D14")
$c.Value = "'4.00"
$c.Style = "Normal"
$ws.Range("E14").Value = "  +0.68%  "

$c = $ws.Range("D15")
$c.Value = "'0.517"
$c.Style = "Normal"
$ws.Range("E15").Value = "  +1.59%  "

$ws.Range("D16").Value = "26.309.98"
$ws.Range("E16").Value = "  +1.37%  "

$c = $ws.Range("D17")
$c.Value = "'62.19"
$c.Style = "Normal"
$ws.Range("E17").Value = "  +3.59%  "

$ws.Range("E18").Value = "  +0.94%  "

$ws.Range("E19").Value = "  -0.01%  "

$c = $ws.Range("D20")
$c.Value = "'201.47"
$c.Style = "Normal"
$ws.Range("E20").Value = "  +1.12%  "

$c = $ws.Range("D22")
$c.Value = "'9.32"
$c.Style = "Normal"
$ws.Range("E22").Value = "  +1.73%  "

$c = $ws.Range("D23")
$c.Value = "'6.04"
$c.Style = "Normal"
$ws.Range("E23").Value = "  +1.23%  "

$c = $ws.Range("D24")
$c.Value = "'1.88"
$c.Style = "Normal"
$ws.Range("E24").Value = "  +1.69%  "

$c = $ws.Range("D25")
$c.Value = "'144.62"
$c.Style = "Normal"
$ws.Range("E25").Value = "  +1.41%  "

$ws.Range("E26").Value = "  -0.02%  "

$ws.Range("E27").Value = "  -1.25%  "

$ws.Range("E28").Value = "  +0.73%  "

$c = $ws.Range("D30")
$c.Value = "'0.0522"
$c.Style = "Normal"
$ws.Range("E30").Value = "  +10.41%  "

$ws.Range("E31").Value = "  +0.96%  "

$ws.Range("E32").Value = "  +2.01%  "

$ws.Range("E34").Value = "  +1.91%  "

$ws.Range("E35").Value = "  +2.37%  "

$ws.Range("D36").Value = "1.179.83"
$ws.Range("E36").Value = "  +5.10%  "

$c = $ws.Range("D37")
$c.Value = "'0.0163"
$c.Style = "Normal"
$ws.Range("E37").Value = "  +0.12%  "

$c = $ws.Range("D38")
$c.Value = "'0.806"
$c.Style = "Normal"
$ws.Range("E38").Value = "  +3.02%  "

$ws.Range("E39").Value = "  -0.04%  "

$ws.Range("E40").Value = "  +0.21%  "

$ws.Range("E41").Value = "  +1.60%  "

$c = $ws.Range("D42")
$c.Value = "'0.787"
$c.Style = "Normal"
$ws.Range("E42").Value = "  +1.46%  "

$ws.Range("E43").Value = "  +4.88%  "

$ws.Range("D44").Value = "1.756.49"
$ws.Range("E44").Value = "  +2.09%  "

$c = $ws.Range("D45")
$c.Value = "'92.60"
$c.Style = "Normal"
$ws.Range("E45").Value = "  +0.81%  "

$ws.Range("B46").Value = "RenderToken"
$ws.Range("C46").Value = "https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr"
$c = $ws.Range("D46")
$c.Value = "'1.53"
$c.Style = "Normal"
$ws.Range("E46").Value = "  +3.46%  "

$ws.Range("B47").Value = "Aave"
$ws.Range("C47").Value = "https://coinranking.com/coin/ixgUfzmLR+aave-aave"
$c = $ws.Range("D47")
$c.Value = "'53.75"
$c.Style = "Normal"
$ws.Range("E47").Value = "  +0.99%  "

$ws.Range("B48").Value = "Cronos"
$ws.Range("C48").Value = "https://coinranking.com/coin/65PHZTpmE55b+cronos-cro"
$c = $ws.Range("D48")
$c.Value = "'0.0508"
$c.Style = "Normal"
$ws.Range("E48").Value = "  +1.09%  "

$ws.Range("B49").Value = "Mantle"
$ws.Range("C49").Value = "https://coinranking.com/coin/BoI4ux0nd+mantle-mnt"
$c = $ws.Range("D49")
$c.Value = "'0.408"
$c.Style = "Normal"
$ws.Range("E49").Value = "  +0.44%  "

$ws.Range("B50").Value = "USDD"
$ws.Range("C50").Value = "https://coinranking.com/coin/z2PZIKQL7+usdd-usdd"
$c = $ws.Range("D50")
$c.Value = "'1.00"
$c.Style = "Normal"
$ws.Range("E50").Value = "  -0.26%  "

$ws.Range("B51").Value = "EnergySwap"
$ws.Range("C51").Value = "https://coinranking.com/coin/SbWqqTui-+energyswap-ens"
$c = $ws.Range("D51")
$c.Value = "'7.28"
$c.Style = "Normal"
$ws.Range("E51").Value = "  +3.06%  "
